# Update NATMI ligand-receptor edge metrics with recomputed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 21.64449733333333
$ws.Range("H2").Value = 64.933492
$ws.Range("I2").Value = 0.1098676276771345
$ws.Range("J2").Value = 0.1098676276771345
$ws.Range("M2").Value = 0.2901893333333334
$ws.Range("N2").Value = 0.870568
$ws.Range("O2").Value = 0.03429389578125064
$ws.Range("P2").Value = 0.03429389578125064
$ws.Range("Q2").Value = 6.281002251495112
$ws.Range("R2").Value = 56.529020263456
$ws.Range("S2").Value = 0.003767788973292899
$ws.Range("T2").Value = 0.003767788973292899

# Row 3
$ws.Range("G3").Value = 21.64449733333333
$ws.Range("H3").Value = 64.933492
$ws.Range("I3").Value = 0.1098676276771345
$ws.Range("J3").Value = 0.1098676276771345
$ws.Range("O3").Value = 0.8402845891331153
$ws.Range("P3").Value = 0.8402845891331153
$ws.Range("Q3").Value = 153.8999660437316
$ws.Range("R3").Value = 1385.099694393584
$ws.Range("S3").Value = 0.09232007438171105
$ws.Range("T3").Value = 0.09232007438171107

# Row 4
$ws.Range("G4").Value = 21.64449733333333
$ws.Range("H4").Value = 64.933492
$ws.Range("I4").Value = 0.1098676276771345
$ws.Range("J4").Value = 0.1098676276771345
$ws.Range("O4").Value = 0.1254215150856341
$ws.Range("P4").Value = 0.1254215150856341
$ws.Range("Q4").Value = 22.97122565670978
$ws.Range("R4").Value = 206.741030910388
$ws.Range("S4").Value = 0.01377976432213056
$ws.Range("T4").Value = 0.01377976432213056

# Row 5
$ws.Range("H5").Value = 88.285005
$ws.Range("I5").Value = 0.1493784449296822
$ws.Range("J5").Value = 0.1493784449296822
$ws.Range("M5").Value = 0.2901893333333334
$ws.Range("N5").Value = 0.870568
$ws.Range("O5").Value = 0.03429389578125064
$ws.Range("P5").Value = 0.03429389578125064
$ws.Range("Q5").Value = 8.53978891476
$ws.Range("R5").Value = 76.85810023284
$ws.Range("S5").Value = 0.005122768822383808
$ws.Range("T5").Value = 0.005122768822383809

# Row 6
$ws.Range("H6").Value = 88.285005
$ws.Range("I6").Value = 0.1493784449296822
$ws.Range("J6").Value = 0.1493784449296822
$ws.Range("O6").Value = 0.8402845891331153
$ws.Range("P6").Value = 0.8402845891331153
$ws.Range("S6").Value = 0.1255204052230817
$ws.Range("T6").Value = 0.1255204052230817

# Row 7
$ws.Range("H7").Value = 88.285005
$ws.Range("I7").Value = 0.1493784449296822
$ws.Range("J7").Value = 0.1493784449296822
$ws.Range("O7").Value = 0.1254215150856341
$ws.Range("P7").Value = 0.1254215150856341
$ws.Range("S7").Value = 0.0187352708842167
$ws.Range("T7").Value = 0.0187352708842167

# Row 8
$ws.Range("I8").Value = 0.7407539273931834
$ws.Range("J8").Value = 0.7407539273931834
$ws.Range("M8").Value = 0.2901893333333334
$ws.Range("N8").Value = 0.870568
$ws.Range("O8").Value = 0.03429389578125064
$ws.Range("P8").Value = 0.03429389578125064
$ws.Range("Q8").Value = 42.34802538408445
$ws.Range("R8").Value = 381.13222845676
$ws.Range("S8").Value = 0.02540333798557393
$ws.Range("T8").Value = 0.02540333798557393

# Row 9
$ws.Range("I9").Value = 0.7407539273931834
$ws.Range("J9").Value = 0.7407539273931834
$ws.Range("O9").Value = 0.8402845891331153
$ws.Range("P9").Value = 0.8402845891331153
$ws.Range("S9").Value = 0.6224441095283226
$ws.Range("T9").Value = 0.6224441095283226

# Row 10
$ws.Range("I10").Value = 0.7407539273931834
$ws.Range("J10").Value = 0.7407539273931834
$ws.Range("O10").Value = 0.1254215150856341
$ws.Range("P10").Value = 0.1254215150856341
$ws.Range("S10").Value = 0.09290647987928687
$ws.Range("T10").Value = 0.09290647987928685
